$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixture")

# The first fixture row (Manchester United v Middlesbrough, FA Cup) is removed;
# every following row shifts up one position.
$ws.Rows.Item(1).Delete()
